$d = $word.ActiveDocument

# 1. "pa deklará e komportashon" -> "pa menshoná e komportashon"
$d.Content.Find.Execute("mayor pa deklará e komportashon", $true, $false, $false, $false, $false,
                         $true, 1, $false, "mayor pa menshoná e komportashon", 2)

# 1b. "palabranan amistoso ora e ta papia" -> "palabranan amabel ora e ta papia"
$d.Content.Find.Execute("palabranan amistoso ora e ta papia", $true, $false, $false, $false, $false,
                         $true, 1, $false, "palabranan amabel ora e ta papia", 2)

# 2. "yud'é bisti mainta" -> "yud'é bisti paña mainta"
$d.Content.Find.Execute("yud’é bisti mainta", $true, $false, $false, $false, $false,
                         $true, 1, $false, "yud’é bisti paña mainta", 2)

# 3. "famia si ta presente" -> "famia si nan ta presente"
$d.Content.Find.Execute("famia si ta presente", $true, $false, $false, $false, $false,
                         $true, 1, $false, "famia si nan ta presente", 2)

# 4. "famia yegá ku bo por konfia den serkania?" -> "famia yegá ku bo por konfia?"
$d.Content.Find.Execute("famia yegá ku bo por konfia den serkania?", $true, $false, $false, $false, $false,
                         $true, 1, $false, "famia yegá ku bo por konfia?", 2)

# 5. "Kòrda e mayor ku e lo ta enfokando riba e yu ku el a skohe pa enfoká riba dje durante e evaluashon."
#    -> "Kòrda e mayor ku e lo enfoká riba e yu ku el a skohe durante di e evaluashon."
$d.Content.Find.Execute("Kòrda e mayor ku e lo ta enfokando riba e yu ku el a skohe pa enfoká riba dje durante e evaluashon.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Kòrda e mayor ku e lo enfoká riba e yu ku el a skohe durante di e evaluashon.", 2)

# 6. "konfrontando pa loke ta trata manehá e komportashon" -> "konfrontando pa manehá e komportashon"
$d.Content.Find.Execute("konfrontando pa loke ta trata manehá e komportashon", $true, $false, $false, $false, $false,
                         $true, 1, $false, "konfrontando pa manehá e komportashon", 2)
